# "changed camera issue image" -- relabel the two callout textboxes on the
# camera-issue diagram slide and nudge their auto-fit widths to match the
# new (slightly different length) text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$EMUS_PER_POINT = 12700

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $text = $shp.TextFrame.TextRange.Text
        if ($text -eq "Real object") {
            $shp.TextFrame.TextRange.Text = "Real Scene"
            $shp.Width = 2266967 / $EMUS_PER_POINT
        } elseif ($text -eq "Camera view") {
            $shp.TextFrame.TextRange.Text = "Camera View"
            $shp.Width = 2646045 / $EMUS_PER_POINT
        }
    }
}
